# Fruta / hortaliza, semanal
# A new weekly price record is added at the top of the data table (row 8),
# pushing all existing records (rows 8-124) down by one row to (9-125).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 8; this shifts rows 8:124 -> 9:125
# and extends the used range / dimension to A1:T125 automatically.
$ws.Rows(8).Insert()

# Populate the newly inserted row 8 with this week's record.
$ws.Range("A8").Value = 10
$ws.Range("B8").Value = "Vega Modelo de Temuco"
$ws.Range("C8").Value = "La Araucanía"
$ws.Range("D8").Value = 45163
$ws.Range("E8").Value = 9
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100108
$ws.Range("H8").Value = "Tropicales y subtropicales"
$ws.Range("I8").Value = 100108007
$ws.Range("J8").Value = "Coco"
$ws.Range("K8").Value = "Sin especificar"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 20
$ws.Range("N8").Value = 36000
$ws.Range("O8").Value = 36000
$ws.Range("P8").Value = 36000
$ws.Range("Q8").Value = '$/malla 20 unidades'
$ws.Range("R8").Value = "Perú"
$ws.Range("S8").Value = 1800
$ws.Range("T8").Value = 20
